# Group 28 (requirements 177-220): populate column I ("Unable to Test" /
# "Success (...)" / "Defect") results for rows 97-224, per the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy formatting (cell style) onto the new I-column cells ---
# style "8" (bordered / Aptos Narrow font) mirrors the existing N137 cell;
# style "4" (bordered, default font) mirrors the existing N154 cell.
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(137, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(138, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(140, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(141, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(143, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(148, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(149, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(152, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(153, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N154").Copy() | Out-Null
$ws.Cells.Item(154, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(156, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(157, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(158, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(159, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(160, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(161, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(163, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(166, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(169, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(170, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(171, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(172, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(173, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(174, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(179, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(182, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(183, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(184, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(185, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(186, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(187, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(188, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(189, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(190, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(191, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(192, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(193, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(194, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(195, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(196, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(197, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(198, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(199, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(200, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(201, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(202, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(203, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(204, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(205, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(206, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(207, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(208, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(209, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(210, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(211, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(212, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(213, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(214, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(216, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(217, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(218, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(219, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(220, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(221, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(222, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(223, 9).PasteSpecial(-4122) | Out-Null
$ws.Range("N137").Copy() | Out-Null
$ws.Cells.Item(224, 9).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 2: write the cell values (reuses/creates shared strings) ---
$ws.Cells.Item(97, 9).Value = "Unable to Test"
$ws.Cells.Item(101, 9).Value = "Unable to Test"
$ws.Cells.Item(102, 9).Value = "Unable to Test"
$ws.Cells.Item(106, 9).Value = "Unable to Test"
$ws.Cells.Item(116, 9).Value = "Unable to Test"
$ws.Cells.Item(120, 9).Value = "Unable to Test"
$ws.Cells.Item(121, 9).Value = "Unable to Test"
$ws.Cells.Item(122, 9).Value = "Unable to Test"
$ws.Cells.Item(123, 9).Value = "Unable to Test"
$ws.Cells.Item(125, 9).Value = "Unable to Test"
$ws.Cells.Item(126, 9).Value = "Unable to Test"
$ws.Cells.Item(129, 9).Value = "Unable to Test"
$ws.Cells.Item(131, 9).Value = "Unable to Test"
$ws.Cells.Item(132, 9).Value = "Unable to Test"
$ws.Cells.Item(133, 9).Value = "Unable to Test"
$ws.Cells.Item(134, 9).Value = "Unable to Test"
$ws.Cells.Item(137, 9).Value = "Unable to Test"
$ws.Cells.Item(138, 9).Value = "Unable to Test"
$ws.Cells.Item(140, 9).Value = "Unable to Test"
$ws.Cells.Item(141, 9).Value = "Unable to Test"
$ws.Cells.Item(143, 9).Value = "Unable to Test"
$ws.Cells.Item(148, 9).Value = "Unable to Test"
$ws.Cells.Item(149, 9).Value = "Unable to Test"
$ws.Cells.Item(152, 9).Value = "Unable to Test"
$ws.Cells.Item(153, 9).Value = "Unable to Test"
$ws.Cells.Item(156, 9).Value = "Unable to Test"
$ws.Cells.Item(157, 9).Value = "Unable to Test"
$ws.Cells.Item(158, 9).Value = "Unable to Test"
$ws.Cells.Item(159, 9).Value = "Unable to Test"
$ws.Cells.Item(160, 9).Value = "Unable to Test"
$ws.Cells.Item(161, 9).Value = "Unable to Test"
$ws.Cells.Item(163, 9).Value = "Unable to Test"
$ws.Cells.Item(166, 9).Value = "Unable to Test"
$ws.Cells.Item(169, 9).Value = "Unable to Test"
$ws.Cells.Item(170, 9).Value = "Unable to Test"
$ws.Cells.Item(171, 9).Value = "Unable to Test"
$ws.Cells.Item(172, 9).Value = "Unable to Test"
$ws.Cells.Item(173, 9).Value = "Unable to Test"
$ws.Cells.Item(174, 9).Value = "Unable to Test"
$ws.Cells.Item(179, 9).Value = "Unable to Test"
$ws.Cells.Item(181, 9).Value = "Success (Test 1: Step 15)"
$ws.Cells.Item(182, 9).Value = "Unable to Test"
$ws.Cells.Item(183, 9).Value = "Unable to Test"
$ws.Cells.Item(184, 9).Value = "Unable to Test"
$ws.Cells.Item(185, 9).Value = "Unable to Test"
$ws.Cells.Item(186, 9).Value = "Unable to Test"
$ws.Cells.Item(187, 9).Value = "Unable to Test"
$ws.Cells.Item(188, 9).Value = "Unable to Test"
$ws.Cells.Item(189, 9).Value = "Unable to Test"
$ws.Cells.Item(190, 9).Value = "Unable to Test"
$ws.Cells.Item(191, 9).Value = "Unable to Test"
$ws.Cells.Item(192, 9).Value = "Unable to Test"
$ws.Cells.Item(193, 9).Value = "Unable to Test"
$ws.Cells.Item(194, 9).Value = "Unable to Test"
$ws.Cells.Item(195, 9).Value = "Unable to Test"
$ws.Cells.Item(196, 9).Value = "Unable to Test"
$ws.Cells.Item(197, 9).Value = "Unable to Test"
$ws.Cells.Item(198, 9).Value = "Unable to Test"
$ws.Cells.Item(199, 9).Value = "Unable to Test"
$ws.Cells.Item(200, 9).Value = "Unable to Test"
$ws.Cells.Item(201, 9).Value = "Unable to Test"
$ws.Cells.Item(202, 9).Value = "Unable to Test"
$ws.Cells.Item(203, 9).Value = "Unable to Test"
$ws.Cells.Item(204, 9).Value = "Unable to Test"
$ws.Cells.Item(205, 9).Value = "Unable to Test"
$ws.Cells.Item(206, 9).Value = "Unable to Test"
$ws.Cells.Item(207, 9).Value = "Unable to Test"
$ws.Cells.Item(208, 9).Value = "Success (Test 1: Step 3)"
$ws.Cells.Item(209, 9).Value = "Success (Test1: Step 3&4)"
$ws.Cells.Item(210, 9).Value = "Success (based off code)"
$ws.Cells.Item(211, 9).Value = "Success (Test 1: Step 5&6)"
$ws.Cells.Item(212, 9).Value = "Success (Test 1: Step 3&4)"
$ws.Cells.Item(213, 9).Value = "Success (Test 1: Step 8)"
$ws.Cells.Item(214, 9).Value = "Success (Test 1: Step 8)"
$ws.Cells.Item(216, 9).Value = "Unable to Test"
$ws.Cells.Item(217, 9).Value = "Unable to Test"
$ws.Cells.Item(218, 9).Value = "Unable to Test"
$ws.Cells.Item(219, 9).Value = "Success (Test 1: Step 3)"
$ws.Cells.Item(220, 9).Value = "Success (Test 1: Step 2&3)"
$ws.Cells.Item(221, 9).Value = "Success (Test 1: Step 4)"
$ws.Cells.Item(222, 9).Value = "Success (Test 1: Step 3)"
$ws.Cells.Item(223, 9).Value = "Defect "
$ws.Cells.Item(224, 9).Value = "Success (Test 1: Step 11)"

# --- Step 3: restore the selection the author left (I95:I224, active I95) ---
$ws.Range("I95:I224").Select()

Write-Output "Group 28 test results written."
